$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (index, distancia, max, min, tempo)
$data = @(
    @(0, 2913.066666666667, 3201, 2602, 0.0550059715906779),
    @(1, 2815.3, 3117, 2493, 0.05445197423299154),
    @(2, 2930.6, 3208, 2688, 0.0589149554570516),
    @(3, 3041.4, 3301, 2665, 0.05837992032368978),
    @(4, 2927.766666666667, 3187, 2708, 0.05862931410471598),
    @(5, 2927.9, 3247, 2527, 0.05870193640391032),
    @(6, 3279.933333333333, 3587, 3063, 0.05870157877604167),
    @(7, 2692.466666666667, 2913, 2383, 0.05611888567606608),
    @(8, 3034.5, 3322, 2496, 0.05596768061319987),
    @(9, 2836.7, 3081, 2449, 0.0557794729868571)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
